# Applies the NATMI re-run update (new "Resolving-Mac" cluster) to the
# Col4a3-Cd47 LR-pairs worksheet: refreshes all recalculated TPM-derived
# metrics for existing sending/target clusters and appends the five new
# rows produced by the newly-added "Resolving-Mac" sending cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2, "ECs", "Col4a3", "Cd47", "ECs", 3, 1, 0.165709, 0.497127, 0.4546154542569759, 0.4546154542569759, 3, 1, 56.756364, 170.269092, 0.157357217290148, 0.157357217290148, 9.405040322076, 84.645362898684, 0.07153702281897428, 0.07153702281897428),
  @(3, "ECs", "Col4a3", "Cd47", "FAPs", 3, 1, 0.165709, 0.497127, 0.4546154542569759, 0.4546154542569759, 3, 1, 71.14312966666667, 213.429389, 0.1972445753159741, 0.1972445753159741, 11.78905687393367, 106.101511865403, 0.08967043220699585, 0.08967043220699585),
  @(4, "ECs", "Col4a3", "Cd47", "Inflammatory-Mac", 3, 1, 0.165709, 0.497127, 0.4546154542569759, 0.4546154542569759, 3, 1, 124.7878343333334, 374.363503, 0.3459747062436438, 0.3459747062436438, 20.67846723954234, 186.106205155881, 0.157285448240378, 0.1572854482403779),
  @(5, "ECs", "Col4a3", "Cd47", "MuSCs", 3, 1, 0.165709, 0.497127, 0.4546154542569759, 0.4546154542569759, 3, 1, 21.495283, 64.485849, 0.05959574714377799, 0.05959574714377799, 3.561961850647, 32.057656655823, 0.0270931476595525, 0.0270931476595525),
  @(6, "ECs", "Col4a3", "Cd47", "Resolving-Mac", 3, 1, 0.165709, 0.497127, 0.4546154542569759, 0.4546154542569759, 3, 1, 86.502237, 259.506711, 0.239827754006456, 0.2398277540064561, 14.334199191033, 129.007792719297, 0.1090294033310753, 0.1090294033310753),
  @(7, "FAPs", "Col4a3", "Cd47", "ECs", 2, 0.6666666666666666, 0.1058106666666667, 0.317432, 0.2902869747080734, 0.2902869747080733, 3, 1, 56.756364, 170.269092, 0.157357217290148, 0.157357217290148, 6.005428712415999, 54.048858411744, 0.04567875055563798, 0.04567875055563798),
  @(8, "FAPs", "Col4a3", "Cd47", "FAPs", 2, 0.6666666666666666, 0.1058106666666667, 0.317432, 0.2902869747080734, 0.2902869747080733, 3, 1, 71.14312966666667, 213.429389, 0.1972445753159741, 0.1972445753159741, 7.527701978783111, 67.74931780904801, 0.05725753104605283, 0.05725753104605283),
  @(9, "FAPs", "Col4a3", "Cd47", "Inflammatory-Mac", 2, 0.6666666666666666, 0.1058106666666667, 0.317432, 0.2902869747080734, 0.2902869747080733, 3, 1, 124.7878343333334, 374.363503, 0.3459747062436438, 0.3459747062436438, 13.20388394269956, 118.834955484296, 0.1004319508009817, 0.1004319508009817),
  @(10, "FAPs", "Col4a3", "Cd47", "MuSCs", 2, 0.6666666666666666, 0.1058106666666667, 0.317432, 0.2902869747080734, 0.2902869747080733, 3, 1, 21.495283, 64.485849, 0.05959574714377799, 0.05959574714377799, 2.274430224418666, 20.469872019768, 0.01729986914383462, 0.01729986914383461),
  @(11, "FAPs", "Col4a3", "Cd47", "Resolving-Mac", 2, 0.6666666666666666, 0.1058106666666667, 0.317432, 0.2902869747080734, 0.2902869747080733, 3, 1, 86.502237, 259.506711, 0.239827754006456, 0.2398277540064561, 9.152859365127998, 82.375734286152, 0.06961887316156615, 0.06961887316156613),
  @(12, "Inflammatory-Mac", "Col4a3", "Cd47", "ECs", 1, 0.3333333333333333, 0.001077, 0.003231, 0.002954702787626279, 0.002954702787626278, 3, 1, 56.756364, 170.269092, 0.157357217290148, 0.157357217290148, 0.061126604028, 0.550139436252, 0.0004649438085803143, 0.0004649438085803142),
  @(13, "Inflammatory-Mac", "Col4a3", "Cd47", "FAPs", 1, 0.3333333333333333, 0.001077, 0.003231, 0.002954702787626279, 0.002954702787626278, 3, 1, 71.14312966666667, 213.429389, 0.1972445753159741, 0.1972445753159741, 0.076621150651, 0.6895903558590001, 0.0005827990965302702, 0.00058279909653027),
  @(14, "Inflammatory-Mac", "Col4a3", "Cd47", "Inflammatory-Mac", 1, 0.3333333333333333, 0.001077, 0.003231, 0.002954702787626279, 0.002954702787626278, 3, 1, 124.7878343333334, 374.363503, 0.3459747062436438, 0.3459747062436438, 0.134396497577, 1.209568478193, 0.001022252428986278, 0.001022252428986277),
  @(15, "Inflammatory-Mac", "Col4a3", "Cd47", "MuSCs", 1, 0.3333333333333333, 0.001077, 0.003231, 0.002954702787626279, 0.002954702787626278, 3, 1, 21.495283, 64.485849, 0.05959574714377799, 0.05959574714377799, 0.023150419791, 0.208353778119, 0.0001760877202163917, 0.0001760877202163917),
  @(16, "Inflammatory-Mac", "Col4a3", "Cd47", "Resolving-Mac", 1, 0.3333333333333333, 0.001077, 0.003231, 0.002954702787626279, 0.002954702787626278, 3, 1, 86.502237, 259.506711, 0.239827754006456, 0.2398277540064561, 0.093162909249, 0.838466183241, 0.0007086197333130252, 0.0007086197333130251),
  @(17, "MuSCs", "Col4a3", "Cd47", "ECs", 3, 1, 0.07261866666666666, 0.217856, 0.1992261623339866, 0.1992261623339865, 3, 1, 56.756364, 170.269092, 0.157357217290148, 0.157357217290148, 4.121571478528, 37.094143306752, 0.03134967451627142, 0.03134967451627142),
  @(18, "MuSCs", "Col4a3", "Cd47", "FAPs", 3, 1, 0.07261866666666666, 0.217856, 0.1992261623339866, 0.1992261623339865, 3, 1, 71.14312966666667, 213.429389, 0.1972445753159741, 0.1972445753159741, 5.166319218887111, 46.496872969984, 0.03929627978139849, 0.03929627978139849),
  @(19, "MuSCs", "Col4a3", "Cd47", "Inflammatory-Mac", 3, 1, 0.07261866666666666, 0.217856, 0.1992261623339866, 0.1992261623339865, 3, 1, 124.7878343333334, 374.363503, 0.3459747062436438, 0.3459747062436438, 9.061926145507556, 81.557335309568, 0.06892721298954951, 0.0689272129895495),
  @(20, "MuSCs", "Col4a3", "Cd47", "MuSCs", 3, 1, 0.07261866666666666, 0.217856, 0.1992261623339866, 0.1992261623339865, 3, 1, 21.495283, 64.485849, 0.05959574714377799, 0.05959574714377799, 1.560958791082667, 14.048629119744, 0.01187303199488153, 0.01187303199488153),
  @(21, "MuSCs", "Col4a3", "Cd47", "Resolving-Mac", 3, 1, 0.07261866666666666, 0.217856, 0.1992261623339866, 0.1992261623339865, 3, 1, 86.502237, 259.506711, 0.239827754006456, 0.2398277540064561, 6.281677114623999, 56.535094031616, 0.04777996305188561, 0.04777996305188561),
  @(22, "Resolving-Mac", "Col4a3", "Cd47", "ECs", 1, 0.3333333333333333, 0.01928833333333333, 0.057865, 0.05291670591333787, 0.05291670591333786, 3, 1, 56.756364, 170.269092, 0.157357217290148, 0.157357217290148, 1.09473566762, 9.85262100858, 0.008326825590683963, 0.008326825590683965),
  @(23, "Resolving-Mac", "Col4a3", "Cd47", "FAPs", 1, 0.3333333333333333, 0.01928833333333333, 0.057865, 0.05291670591333787, 0.05291670591333786, 3, 1, 71.14312966666667, 213.429389, 0.1972445753159741, 0.1972445753159741, 1.372232399387222, 12.350091594485, 0.01043753318499662, 0.01043753318499662),
  @(24, "Resolving-Mac", "Col4a3", "Cd47", "Inflammatory-Mac", 1, 0.3333333333333333, 0.01928833333333333, 0.057865, 0.05291670591333787, 0.05291670591333786, 3, 1, 124.7878343333334, 374.363503, 0.3459747062436438, 0.3459747062436438, 2.406949344566112, 21.662544101095, 0.01830784178374836, 0.01830784178374836),
  @(25, "Resolving-Mac", "Col4a3", "Cd47", "MuSCs", 1, 0.3333333333333333, 0.01928833333333333, 0.057865, 0.05291670591333787, 0.05291670591333786, 3, 1, 21.495283, 64.485849, 0.05959574714377799, 0.05959574714377799, 0.4146081835983333, 3.731473652385, 0.003153610625292945, 0.003153610625292945),
  @(26, "Resolving-Mac", "Col4a3", "Cd47", "Resolving-Mac", 1, 0.3333333333333333, 0.01928833333333333, 0.057865, 0.05291670591333787, 0.05291670591333786, 3, 1, 86.502237, 259.506711, 0.239827754006456, 0.2398277540064561, 1.668483981335, 15.016355832015, 0.01269089472861597, 0.01269089472861597)
)

foreach ($row in $data) {
  $r = $row[0]
  for ($c = 1; $c -lt $row.Count; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c]
  }
}

